# Applies the "new api calls and some edits" change:
#  - splits editUser into editNameUser + editPassUser (new sections),
#    while keeping the original editUser section intact afterwards
#  - marks getOneContact as DEPRECATED

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the editNameUser + editPassUser sections right before the
#    existing "editUser:" heading paragraph.
# ---------------------------------------------------------------------------

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "editUser:`r") {
        $targetIndex = $i
        break
    }
}

$editUserHeading = $d.Paragraphs.Item($targetIndex)
$editUserHeading.Range.InsertParagraphBefore()
$insertionPara = $d.Paragraphs.Item($targetIndex)

$newSectionsXml = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>editNameUser</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">required input: </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>user_id</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>full_name</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, and </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>new_name</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">output: </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>user_id</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>full_name</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>, status, and message</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>editPassUser</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">required input: </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>user_id</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, password, </w:t></w:r>
  <w:r><w:t xml:space="preserve">and </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>new_pass</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>
  <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">output: </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>user_id</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>,</w:t></w:r>
  <w:bookmarkStart w:id="100" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="100"/>
  <w:r><w:t xml:space="preserve"> status, and message</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPara.Range.InsertXML($newSectionsXml)

# ---------------------------------------------------------------------------
# 2) Mark getOneContact as DEPRECATED.
# ---------------------------------------------------------------------------

$getOneContactHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "getOneContact:`r") {
        $getOneContactHeading = $para
        break
    }
}

$tailRange = $getOneContactHeading.Range.Duplicate
$tailRange.Collapse(0)
$tailRange.MoveEnd(1, -1)
$tailRange.Text = " DEPRECATED"
$tailRange.Font.Bold = $true

Write-Output "done"
